$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Code Review 3" column (D): each of the 5 students scored 20/20
$ws.Range("D3:D7").Value = 20

# Total row (row 8) for the new column
$ws.Range("D8").Value = 100

# Leave the selection where the author left it when saving
$ws.Range("E3").Select()
